# Auto-generated edit script: update crypto price/volume table cells
# per commit "Updated symbol list on Mon Jan  2 14:43:43 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '246.28'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.81%'
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '29.82'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '9.67%'
$ws.Range('E3').ClearFormats()
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.177'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.42%'
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05701'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '0.65%'
$ws.Range('E5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.591'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '1.89%'
$ws.Range('E6').ClearFormats()
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.050'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '1.75%'
$ws.Range('E7').ClearFormats()
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8584'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '4.40%'
$ws.Range('E8').ClearFormats()
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8751'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '3.72%'
$ws.Range('E9').ClearFormats()
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1367'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '3.21%'
$ws.Range('E10').ClearFormats()
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07095'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '2.45%'
$ws.Range('E11').ClearFormats()
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.02865'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-0.69%'
$ws.Range('E12').ClearFormats()
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09382'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.08%'
$ws.Range('E13').ClearFormats()
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001530'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '1.04%'
$ws.Range('E14').ClearFormats()
$ws.Range('B15').Value = 'CoinExToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.04156'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.89%'
$ws.Range('E15').ClearFormats()
$ws.Range('B16').Value = 'One'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0006011'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-94.00%'
$ws.Range('E16').ClearFormats()
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.006001'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-3.14%'
$ws.Range('E17').ClearFormats()
$ws.Range('B18').Value = 'UpBots'
$ws.Range('C18').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.007491'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '5,107.63%'
$ws.Range('E18').ClearFormats()
$ws.Range('B19').Value = 'LEO'
$ws.Range('C19').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.481'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-0.92%'
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.264'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-2.04%'
$ws.Range('E20').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '3.60%'
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.1300'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '0.72%'
$ws.Range('E23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.905'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-18.29%'
$ws.Range('E24').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '0.46%'
$ws.Range('E25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.005086'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '31.41%'
$ws.Range('E26').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.001222'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '0.34%'
$ws.Range('E27').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0001211'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '23.54%'
$ws.Range('E28').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03746'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '1.97%'
$ws.Range('E40').ClearFormats()
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.005783'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-4.41%'
$ws.Range('E41').ClearFormats()
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1071'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '1.75%'
$ws.Range('E42').ClearFormats()
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002101'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-8.66%'
$ws.Range('E43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.009415'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '1.21%'
$ws.Range('E44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005111'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-4.01%'
$ws.Range('E45').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.04%'
$ws.Range('E46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.07101'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-29.69%'
$ws.Range('E47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.002674'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '4.17%'
$ws.Range('E48').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.04%'
$ws.Range('E49').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.04%'
$ws.Range('E50').ClearFormats()
